# Aragon hospitalized/UCI COVID-19 dataset: append two new report dates
# (2020-08-01 => 44044 and 2020-08-03 => 44046) to "Hoja1", one row per
# hospital, replicating the existing table layout/formatting exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last populated row before this edit, and the first new row to write.
$templateRow = 2192
$startRow = 2193

# Columns: A=fecha, B=hospital, C=camas_ocupadas_total, D=camas_uci_ocupadas,
# E=municipio, F=provincia, G=codigo_ine, H=observaciones.
# Rows with no admissions reported for a column (C/D) simply omit that key,
# matching the sparse source rows (e.g. hospitals with 0 UCI beds tracked).
$newRows = @(
    @{ A=44044; B="Hospital Universitario Miguel Servet"; C=110; D=12; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Clínico Universitario"; C=104; D=9; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Royo Villanova"; C=29; D=2; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Nuestra Señora de Gracia"; C=7; D=1; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital General de la Defensa"; C=12; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Obispo Polanco"; C=29; E="Teruel"; F="Teruel"; G=44216; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital de Alcañiz"; C=8; E="Alcañiz"; F="Teruel"; G=44216; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital de Barbastro"; C=20; D=4; E="Barbastro"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital San Jorge"; C=19; D=2; E="Huesca"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Sagrado Corazón"; E="Huesca"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Ernest Lluch"; C=3; E="Calatayud"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital San José"; E="Teruel"; F="Teruel"; G=44216; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Ejea – Cinco Villas"; C=3; E="Ejea de los Caballeros"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="MAZ"; C=3; D=1; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Viamed Montecanal"; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Clínica Montpellier"; C=5; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital Quirón"; C=6; D=1; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Hospital San Juan de Dios de Zaragoza"; C=13; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Clínica Viamed Santiago"; E="Huesca"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44044; B="Clínica El Pilar"; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Universitario Miguel Servet"; C=137; D=15; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Clínico Universitario"; C=135; D=9; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Royo Villanova"; C=34; D=3; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Nuestra Señora de Gracia"; C=14; D=1; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital General de la Defensa"; C=19; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Obispo Polanco"; C=39; E="Teruel"; F="Teruel"; G=44216; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital de Alcañiz"; C=8; E="Alcañiz"; F="Teruel"; G=44216; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital de Barbastro"; C=25; D=4; E="Barbastro"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital San Jorge"; C=24; D=3; E="Huesca"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Sagrado Corazón"; E="Huesca"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Ernest Lluch"; C=4; E="Calatayud"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital San José"; C=0; E="Teruel"; F="Teruel"; G=44216; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Ejea – Cinco Villas"; C=3; E="Ejea de los Caballeros"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="MAZ"; C=4; D=1; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Viamed Montecanal"; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Clínica Montpellier"; C=6; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital Quirón"; C=6; D=3; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Hospital San Juan de Dios de Zaragoza"; C=15; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Clínica Viamed Santiago"; C=1; E="Huesca"; F="Huesca"; G=22125; H="Fuente Aragón Hoy" },
    @{ A=44046; B="Clínica El Pilar"; E="Zaragoza"; F="Zaragoza"; G=50297; H="Fuente Aragón Hoy" }
)

$cols = @("A","B","C","D","E","F","G","H")

$r = $startRow
foreach ($rowData in $newRows) {
    foreach ($col in $cols) {
        if ($rowData.ContainsKey($col)) {
            $cellRef = "$col$r"
            $ws.Range($cellRef).Value = $rowData[$col]
            # Copy the template row's formatting (fill/border/number format)
            # onto the new cell so banding/date format matches the rest of
            # the sheet, without materializing empty cells for blank columns.
            $ws.Range("$col$templateRow").Copy()
            $ws.Range($cellRef).PasteSpecial(-4122)
        }
    }
    $r++
}

$lastRow = $startRow + $newRows.Count - 1
[void]$ws.Range("A2214:A$lastRow").Select()
